$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force the date column to Text format before writing literal date strings,
# so Excel does not auto-convert them into date serial numbers.
$ws.Range("A2:A91").NumberFormat = "@"

$ws.Range("A2").Value = "2025-10-14"
$ws.Range("C2").Value = 41
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("C3").Value = 49
$ws.Range("A4").Value = "2025-10-16"
$ws.Range("C4").Value = 50
$ws.Range("A5").Value = "2025-10-17"
$ws.Range("C5").Value = 59
$ws.Range("A6").Value = "2025-10-18"
$ws.Range("C6").Value = 63
$ws.Range("A7").Value = "2025-10-19"
$ws.Range("C7").Value = 66
$ws.Range("A8").Value = "2025-10-20"
$ws.Range("C8").Value = 72
$ws.Range("A9").Value = "2025-10-21"
$ws.Range("C9").Value = 81
$ws.Range("A10").Value = "2025-10-22"
$ws.Range("C10").Value = 81
$ws.Range("A11").Value = "2025-10-23"
$ws.Range("C11").Value = 83
$ws.Range("A12").Value = "2025-10-24"
$ws.Range("C12").Value = 84
$ws.Range("A13").Value = "2025-10-25"
$ws.Range("C13").Value = 85
$ws.Range("A14").Value = "2025-10-26"
$ws.Range("C14").Value = 90
$ws.Range("A15").Value = "2025-10-27"
$ws.Range("C15").Value = 83
$ws.Range("A16").Value = "2025-10-28"
$ws.Range("C16").Value = 90
$ws.Range("A17").Value = "2025-10-29"
$ws.Range("C17").Value = 93
$ws.Range("A18").Value = "2025-10-30"
$ws.Range("C18").Value = 92
$ws.Range("A19").Value = "2025-10-31"
$ws.Range("C19").Value = 101
$ws.Range("A20").Value = "2025-11-01"
$ws.Range("C20").Value = 112
$ws.Range("A21").Value = "2025-11-02"
$ws.Range("C21").Value = 115
$ws.Range("A22").Value = "2025-11-03"
$ws.Range("C22").Value = 108
$ws.Range("A23").Value = "2025-11-04"
$ws.Range("C23").Value = 105
$ws.Range("A24").Value = "2025-11-05"
$ws.Range("C24").Value = 101
$ws.Range("A25").Value = "2025-11-06"
$ws.Range("C25").Value = 95
$ws.Range("A26").Value = "2025-11-07"
$ws.Range("C26").Value = 87
$ws.Range("A27").Value = "2025-11-08"
$ws.Range("C27").Value = 82
$ws.Range("A28").Value = "2025-11-09"
$ws.Range("C28").Value = 76
$ws.Range("A29").Value = "2025-11-10"
$ws.Range("C29").Value = 54
$ws.Range("A30").Value = "2025-11-11"
$ws.Range("C30").Value = 47
$ws.Range("A31").Value = "2025-11-12"
$ws.Range("C31").Value = 43
$ws.Range("A32").Value = "2025-11-13"
$ws.Range("C32").Value = 41
$ws.Range("A33").Value = "2025-11-14"
$ws.Range("C33").Value = 38
$ws.Range("A34").Value = "2025-11-15"
$ws.Range("C34").Value = 35
$ws.Range("A35").Value = "2025-11-16"
$ws.Range("C35").Value = 31
$ws.Range("A36").Value = "2025-11-17"
$ws.Range("C36").Value = 29
$ws.Range("A37").Value = "2025-11-18"
$ws.Range("C37").Value = 26
$ws.Range("A38").Value = "2025-11-19"
$ws.Range("C38").Value = 26
$ws.Range("A39").Value = "2025-11-20"
$ws.Range("C39").Value = 25
$ws.Range("A40").Value = "2025-11-21"
$ws.Range("C40").Value = 26
$ws.Range("A41").Value = "2025-11-22"
$ws.Range("C41").Value = 26
$ws.Range("A42").Value = "2025-11-23"
$ws.Range("C42").Value = 25
$ws.Range("A43").Value = "2025-11-24"
$ws.Range("C43").Value = 25
$ws.Range("A44").Value = "2025-11-25"
$ws.Range("C44").Value = 27
$ws.Range("A45").Value = "2025-11-26"
$ws.Range("C45").Value = 27
$ws.Range("A46").Value = "2025-11-27"
$ws.Range("C46").Value = 27
$ws.Range("A47").Value = "2025-11-28"
$ws.Range("C47").Value = 27
$ws.Range("A48").Value = "2025-11-29"
$ws.Range("C48").Value = 27
$ws.Range("A49").Value = "2025-11-30"
$ws.Range("C49").Value = 27
$ws.Range("A50").Value = "2025-12-01"
$ws.Range("C50").Value = 27
$ws.Range("A51").Value = "2025-12-02"
$ws.Range("C51").Value = 27
$ws.Range("A52").Value = "2025-12-03"
$ws.Range("C52").Value = 26
$ws.Range("A53").Value = "2025-12-04"
$ws.Range("C53").Value = 25
$ws.Range("A54").Value = "2025-12-05"
$ws.Range("C54").Value = 25
$ws.Range("A55").Value = "2025-12-06"
$ws.Range("C55").Value = 25
$ws.Range("A56").Value = "2025-12-07"
$ws.Range("C56").Value = 26
$ws.Range("A57").Value = "2025-12-08"
$ws.Range("C57").Value = 26
$ws.Range("A58").Value = "2025-12-09"
$ws.Range("C58").Value = 27
$ws.Range("A59").Value = "2025-12-10"
$ws.Range("C59").Value = 29
$ws.Range("A60").Value = "2025-12-11"
$ws.Range("C60").Value = 29
$ws.Range("A61").Value = "2025-12-12"
$ws.Range("C61").Value = 30
$ws.Range("A62").Value = "2025-12-13"
$ws.Range("C62").Value = 30
$ws.Range("A63").Value = "2025-12-14"
$ws.Range("C63").Value = 31
$ws.Range("A64").Value = "2025-12-15"
$ws.Range("C64").Value = 31
$ws.Range("A65").Value = "2025-12-16"
$ws.Range("C65").Value = 32
$ws.Range("A66").Value = "2025-12-17"
$ws.Range("C66").Value = 31
$ws.Range("A67").Value = "2025-12-18"
$ws.Range("C67").Value = 31
$ws.Range("A68").Value = "2025-12-19"
$ws.Range("C68").Value = 32
$ws.Range("A69").Value = "2025-12-20"
$ws.Range("C69").Value = 32
$ws.Range("A70").Value = "2025-12-21"
$ws.Range("C70").Value = 32
$ws.Range("A71").Value = "2025-12-22"
$ws.Range("C71").Value = 32
$ws.Range("A72").Value = "2025-12-23"
$ws.Range("C72").Value = 30
$ws.Range("A73").Value = "2025-12-24"
$ws.Range("C73").Value = 31
$ws.Range("A74").Value = "2025-12-25"
$ws.Range("C74").Value = 32
$ws.Range("A75").Value = "2025-12-26"
$ws.Range("C75").Value = 32
$ws.Range("A76").Value = "2025-12-27"
$ws.Range("C76").Value = 28
$ws.Range("A77").Value = "2025-12-28"
$ws.Range("C77").Value = 28
$ws.Range("A78").Value = "2025-12-29"
$ws.Range("C78").Value = 28
$ws.Range("A79").Value = "2025-12-30"
$ws.Range("C79").Value = 28
$ws.Range("A80").Value = "2025-12-31"
$ws.Range("C80").Value = 30
$ws.Range("A81").Value = "2026-01-01"
$ws.Range("C81").Value = 29
$ws.Range("A82").Value = "2026-01-02"
$ws.Range("C82").Value = 28
$ws.Range("A83").Value = "2026-01-03"
$ws.Range("C83").Value = 28
$ws.Range("A84").Value = "2026-01-04"
$ws.Range("C84").Value = 27
$ws.Range("A85").Value = "2026-01-05"
$ws.Range("C85").Value = 27
$ws.Range("A86").Value = "2026-01-06"
$ws.Range("C86").Value = 27
$ws.Range("A87").Value = "2026-01-07"
$ws.Range("C87").Value = 27
$ws.Range("A88").Value = "2026-01-08"
$ws.Range("C88").Value = 27
$ws.Range("A89").Value = "2026-01-09"
$ws.Range("C89").Value = 27
$ws.Range("A90").Value = "2026-01-10"
$ws.Range("C90").Value = 26
$ws.Range("A91").Value = "2026-01-11"
$ws.Range("C91").Value = 26
$ws.Range("B91").Value = 0.0

# Restore original (default) cell formatting on the date column by
# pasting formats only from an already-default-styled cell.
$ws.Range("B1").Copy()
$ws.Range("A2:A91").PasteSpecial(-4122)
